# Cut Sheet Express - replace shared-string placeholder letters (q,w,e,r,t,y,u,i,f)
# in row 2/3 with real sequential numeric data, extend the sequence through
# columns S:T, clear out row 4 entirely, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: J2:T2 -> 10..20 (numeric, replacing old shared-string letters, plus new S2/T2)
$row2 = New-Object 'object[,]' 1,11
for ($i = 0; $i -lt 11; $i++) { $row2[0,$i] = 10 + $i }
$ws.Range("J2:T2").Value = $row2

# Row 3: A3:T3 -> 21..40 (numeric, replacing old 1-9 and shared-string letters, plus new S3/T3)
$row3 = New-Object 'object[,]' 1,20
for ($i = 0; $i -lt 20; $i++) { $row3[0,$i] = 21 + $i }
$ws.Range("A3:T3").Value = $row3

# Row 4: clear out all the old placeholder values entirely
$null = $ws.Range("A4:R4").ClearContents()

# Update the active selection from H8 to H12
$null = $ws.Range("H12").Select()
